$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values; prefix with an apostrophe so Excel keeps them
# as literal text instead of re-parsing/trimming them as numbers.
$ws.Range("D2").Value = "'27.039.03"
$ws.Range("D3").Value = "'1.891.91"
$ws.Range("D5").Value = "'306.70"
$ws.Range("D6").Value = "'1.003"
$ws.Range("D7").Value = "'0.5229"
$ws.Range("D8").Value = "'0.3757"
$ws.Range("D10").Value = "'21.07"
$ws.Range("D11").Value = "'0.8979"
$ws.Range("D12").Value = "'0.08160"
$ws.Range("D13").Value = "'1.938.02"
$ws.Range("D15").Value = "'5.297"
$ws.Range("D20").Value = "'27.072.28"
$ws.Range("D23").Value = "'6.407"
$ws.Range("D24").Value = "'148.68"
$ws.Range("D25").Value = "'2.286"
$ws.Range("D26").Value = "'18.17"
$ws.Range("D27").Value = "'1.734"
$ws.Range("D28").Value = "'115.01"
$ws.Range("D29").Value = "'4.778"
$ws.Range("D30").Value = "'4.847"
$ws.Range("D31").Value = "'0.09221"
$ws.Range("D32").Value = "'0.05029"
$ws.Range("D33").Value = "'0.7884"
$ws.Range("D35").Value = "'3.422"
$ws.Range("D36").Value = "'2.977"
$ws.Range("D37").Value = "'2.598"
$ws.Range("D38").Value = "'0.5707"
$ws.Range("D41").Value = "'9.031"
$ws.Range("D42").Value = "'6.544"
$ws.Range("D43").Value = "'116.13"
$ws.Range("D45").Value = "'0.4858"
$ws.Range("D46").Value = "'1.003"
$ws.Range("D48").Value = "'1.622"
$ws.Range("D49").Value = "'38.16"
$ws.Range("D50").Value = "'63.46"
$ws.Range("D51").Value = "'0.05929"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  +5.48%  "
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  -0.05%  "
